# Apply the 2024-02-23 strategy sheet edit:
#  - update the set_voltage (column G) readings for rows 3-26
#  - move the active selection/viewport to reflect where the user was
#    last working (topLeftCell A6, active cell G20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = 53
    4  = 51
    5  = 51
    6  = 51
    7  = 51
    8  = 51
    9  = 51
    10 = 51
    11 = 51
    12 = 51
    13 = 51
    14 = 51
    15 = 51
    16 = 51
    17 = 51
    18 = 51
    19 = 55
    20 = 51
    21 = 51
    22 = 51
    23 = 51
    24 = 51
    25 = 51
    26 = 51
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}

# Scroll the view so row 6 is at the top, then move/select the active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("G20").Select()

Write-Output "applied strategy sheet update"
